$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 369, shifting existing rows 369:498 down to 370:499.
$ws.Rows.Item(369).Insert()

# The newly inserted row 369 is blank; copy the (now shifted-down) row 370's
# formatting/values as a baseline, then overwrite the fields that differ for
# the new record.
$ws.Rows.Item(370).Copy()
$ws.Rows.Item(369).PasteSpecial(-4104) | Out-Null   # xlPasteAll
$excel.CutCopyMode = 0

# Populate the new record's values (row 369)
$ws.Range("D369").Value = 44809
$ws.Range("J369").Value = 5500
$ws.Range("K369").Value = 1200
$ws.Range("L369").Value = 1300
$ws.Range("M369").Value = 1245
$ws.Range("P369").Value = 1245
